# Update "paises.xlsx" (countries / provincias Spain) data snapshot.
#
# 1) Refresh the "last updated" timestamp.
# 2) Refresh Australia's stats (row 53).
# 3) "Butan" (Bhutan) case counts increased, moving it up the
#    (descending, sorted-by-total-cases) ranking from its old slot
#    (just before "Montserrat") to a new slot (just before "Burundi").
#    That shifts Burundi / Mauritania / San Cristobal y Nieves /
#    Islas Malvinas / Santa Sede / Islas Turcas y Caicos / Comoras
#    down by one row each; "Montserrat" and everything after keeps
#    its original row/data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) timestamp -----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 05:05"

# --- 2) Australia (row 53): Total, Nuevos, Activos, Recuperados, Criticos
$ws.Range("B53").Value = 6989
$ws.Range("C53").Value = 9
$ws.Range("D53").Value = 6297
$ws.Range("E53").Value = 594
$ws.Range("F53").Value = 18

# --- 3) Butan jumps up the ranking, rows 200-207 shift -----------------
# row 200: Butan (new/updated figures)
$ws.Range("A200").Value = "Butan"
$ws.Range("B200").Value = 15
$ws.Range("C200").Value = 4
$ws.Range("D200").Value = 5
$ws.Range("E200").Value = 10
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 0

# row 201: Burundi (was row 200's data)
$ws.Range("A201").Value = "Burundi"
$ws.Range("B201").Value = 15
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 7
$ws.Range("E201").Value = 7
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 1

# row 202: Mauritania (was row 201's data)
$ws.Range("A202").Value = "Mauritania"
$ws.Range("B202").Value = 15
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 6
$ws.Range("E202").Value = 7
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 2

# row 203: San Cristobal y Nieves (was row 202's data)
$ws.Range("A203").Value = "San Cristobal y Nieves"
$ws.Range("B203").Value = 15
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 14
$ws.Range("E203").Value = 1
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

# row 204: Islas Malvinas (was row 203's data)
$ws.Range("A204").Value = "Islas Malvinas"
$ws.Range("B204").Value = 13
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 13
$ws.Range("E204").Value = 0
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 0

# row 205: Santa Sede (was row 204's data)
$ws.Range("A205").Value = "Santa Sede"
$ws.Range("B205").Value = 12
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 2
$ws.Range("E205").Value = 10
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0

# row 206: Islas Turcas y Caicos (was row 205's data)
$ws.Range("A206").Value = "Islas Turcas y Caicos"
$ws.Range("B206").Value = 12
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 10
$ws.Range("E206").Value = 1
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 1

# row 207: Comoras (was row 206's data)
$ws.Range("A207").Value = "Comoras"
$ws.Range("B207").Value = 11
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 0
$ws.Range("E207").Value = 10
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 1

# row 208 (Montserrat) and beyond are unaffected - old Butan row was
# removed from in front of it and the freed slot absorbed by the shift
# above, so its own row/data stays exactly as it was.
